$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "actual hours worked" (Фактические часы работы) block for the
# second and third day-columns of the first week (row 6 = start time / duration,
# row 7 = end time), mirroring the values already present for the first day.
$ws.Range("E6").Value = 0.625
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 0.46875
$ws.Range("H6").Value = 8.5

$ws.Range("E7").Value = 0.83333333333333337
$ws.Range("G7").Value = 0.82291666666666663

# Reflect the final active selection recorded in the workbook after the edit.
$ws.Range("I6").Select()
